# Generate Report for Handoff
#
# The localization run finished handoff generation: every language's
# status flips from "Handed back: in sync with en-US" to
# "Ready for handoff", and the handoff/generation timestamps on the
# Overview and zh-cn sheets are refreshed to the moment the report was
# (re)generated. Because the new status text is shorter than the old
# one, Excel auto-shrinks the Status-ish columns that held it.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# zh-cn / de-de status columns (E2, F2) and the
# "Latest HO Xliff Generate Date" column (G2).
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-04 21:06:38"

# --- zh-cn sheet ------------------------------------------------------
# Status (C2) and Latest Handoff Datetime (H2).
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-04 21:06:32"

# --- de-de sheet --------------------------------------------------
# Status (C2) only - its Latest Handoff Datetime cell shares the same
# shared string as the Overview's generate date, so it updates along
# with it automatically.
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Column widths ----------------------------------------------------
# The shorter "Ready for handoff" text lets Excel auto-fit these
# columns narrower than the old "Handed back: in sync with en-US".
$wsOverview.Range("E1").ColumnWidth = 16.333333333333336
$wsOverview.Range("F1").ColumnWidth = 16.333333333333336
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333336
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333336
